$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the tab-bar/horizontal-scrollbar split ratio (cosmetic window setting).
$win = $excel.Windows.Item(1)
$win.TabRatio = 0.5

# The "brand" value ("nike") in cell AE2 was incorrect/unwanted for this
# downloadable-product sample row, so select it and clear its contents.
$ws.Range("AE2").Select()
$ws.Range("AE2").ClearContents()

$wb.Save()
